$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1333
$ws.Range("I98").Value = 1499.5
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1499.5
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -1.5
$ws.Range("N98").Value = -3996
$ws.Range("H106").Value = 7999
$ws.Range("I106").Value = 7998.5
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 7998.5
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -7367.5
$ws.Range("N106").Value = -9262
$ws.Range("H107").Value = 591
$ws.Range("I107").Value = 624.4
$ws.Range("K107").Value = 624.4
$ws.Range("M107").Value = 1295.6
$ws.Range("H112").Value = 1618.5652
$ws.Range("J112").Value = 1790.8948
$ws.Range("L112").Value = 5372.6844
$ws.Range("N112").Value = -7588.6844
$ws.Range("H116").Value = 6978.3
$ws.Range("I116").Value = 5901
$ws.Range("K116").Value = 5901
$ws.Range("M116").Value = -2459
$ws.Range("H121").Value = 1948.3572
$ws.Range("J121").Value = 1948.3572
$ws.Range("L121").Value = 5845.071599999999
$ws.Range("N121").Value = -9339.071599999999
$ws.Range("H122").Value = 1333
$ws.Range("I122").Value = 1499.5
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4498.5
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2048.5
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 1350.1333
$ws.Range("I132").Value = 1350.1333
$ws.Range("K132").Value = 4050.3999
$ws.Range("M132").Value = -1520.3999
$ws.Range("H135").Value = 2195.6667
$ws.Range("I135").Value = 2156.8
$ws.Range("K135").Value = 19411.2
$ws.Range("M135").Value = -16876.2
$ws.Range("H138").Value = 2880.894
$ws.Range("I138").Value = 2297.8
$ws.Range("J138").Value = 2985.0178
$ws.Range("K138").Value = 6893.400000000001
$ws.Range("L138").Value = 8955.053400000001
$ws.Range("M138").Value = -1753.400000000001
$ws.Range("N138").Value = -19235.0534
$ws.Range("H141").Value = 5506.8184
$ws.Range("I141").Value = 5506.8184
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 16520.4552
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -11340.4552
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1677.5385
$ws.Range("I2").Value = 1742.3334
$ws.Range("K2").Value = 1742.3334
$ws.Range("M2").Value = -1629.3334
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H116").Value = 1677.5385
$ws.Range("I116").Value = 1742.3334
$ws.Range("K116").Value = 1742.3334
$ws.Range("M116").Value = 551.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1677.5385
$ws.Range("I3").Value = 1742.3334
$ws.Range("K3").Value = 1742.3334
$ws.Range("M3").Value = -1628.3334
$ws.Range("H86").Value = 1672.7273
$ws.Range("I86").Value = 1912.5
$ws.Range("K86").Value = 1912.5
$ws.Range("M86").Value = -789.5
$ws.Range("H89").Value = 1672.7273
$ws.Range("I89").Value = 1912.5
$ws.Range("K89").Value = 9562.5
$ws.Range("M89").Value = -3946.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 17804.783
$ws.Range("I86").Value = 7997.636
$ws.Range("J86").Value = 26794.666
$ws.Range("K86").Value = 7997.636
$ws.Range("L86").Value = 26794.666
$ws.Range("M86").Value = -6874.636
$ws.Range("N86").Value = -29040.666
$ws.Range("H89").Value = 17804.783
$ws.Range("I89").Value = 7997.636
$ws.Range("J89").Value = 26794.666
$ws.Range("K89").Value = 39988.18
$ws.Range("L89").Value = 133973.33
$ws.Range("M89").Value = -34372.18
$ws.Range("N89").Value = -145205.33
$ws.Range("H107").Value = 2531.5833
$ws.Range("J107").Value = 2932.4443
$ws.Range("L107").Value = 2932.4443
$ws.Range("N107").Value = -6772.4443
$ws.Range("H122").Value = 3302
$ws.Range("I122").Value = 2666.3333
$ws.Range("J122").Value = 3937.6667
$ws.Range("K122").Value = 7998.999899999999
$ws.Range("L122").Value = 11813.0001
$ws.Range("M122").Value = -5548.999899999999
$ws.Range("N122").Value = -16713.0001
$ws.Range("H132").Value = 2691.7917
$ws.Range("I132").Value = 2069.4707
$ws.Range("J132").Value = 4203.143
$ws.Range("K132").Value = 6208.4121
$ws.Range("L132").Value = 12609.429
$ws.Range("M132").Value = -3678.4121
$ws.Range("N132").Value = -17669.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 105.5
$ws.Range("I6").Value = 105.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 316.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -203.5
$ws.Range("N6").ClearContents()
$ws.Range("H17").Value = 406.92856
$ws.Range("I17").Value = 399.76923
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 1199.30769
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -1030.30769
$ws.Range("N17").Value = -1838
$ws.Range("H107").Value = 213
$ws.Range("I107").Value = 100
$ws.Range("J107").Value = 235.6
$ws.Range("K107").Value = 300
$ws.Range("L107").Value = 706.8
$ws.Range("M107").Value = 1620
$ws.Range("N107").Value = -4546.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1637
$ws.Range("I102").Value = 1462
$ws.Range("K102").Value = 1462
$ws.Range("M102").Value = 160
$ws.Range("H113").Value = 1949.5
$ws.Range("I113").Value = 1949.5
$ws.Range("K113").Value = 1949.5
$ws.Range("M113").Value = 220.5
$ws.Range("H122").Value = 1679.8
$ws.Range("I122").Value = 1224.75
$ws.Range("K122").Value = 3674.25
$ws.Range("M122").Value = -1224.25
$ws.Range("H132").Value = 2672
$ws.Range("I132").Value = 1740.2
$ws.Range("J132").Value = 3836.75
$ws.Range("K132").Value = 5220.6
$ws.Range("L132").Value = 11510.25
$ws.Range("M132").Value = -2690.6
$ws.Range("N132").Value = -16570.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H61").Value = 2833
$ws.Range("I61").Value = 1999.5
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 1999.5
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -1797.5
$ws.Range("N61").Value = -4904
$ws.Range("H113").Value = 2833
$ws.Range("I113").Value = 1999.5
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 1999.5
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 170.5
$ws.Range("N113").Value = -8840
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 4127.7144
$ws.Range("I136").Value = 4127.7144
$ws.Range("K136").Value = 12383.1432
$ws.Range("M136").Value = -9833.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1693.6666
$ws.Range("I113").Value = 1650
$ws.Range("K113").Value = 4950
$ws.Range("M113").Value = -2780
$ws.Range("H126").Value = 1999.75
$ws.Range("I126").Value = 1999.75
$ws.Range("K126").Value = 5999.25
$ws.Range("M126").Value = -3529.25
